# Re-theme the deck: apply the built-in "Office Theme" colour scheme
# (the swap that PowerPoint performs when a new theme is chosen from the
# Design gallery). This updates the colour scheme that backs every slide
# (ppt/theme/theme1.xml), replacing the "Integral" palette with the
# default "Office" palette.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colours, in ThemeColorScheme item order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$tcs.Item(1).RGB  = HexToRgb "000000"   # dk1
$tcs.Item(2).RGB  = HexToRgb "FFFFFF"   # lt1
$tcs.Item(3).RGB  = HexToRgb "44546A"   # dk2
$tcs.Item(4).RGB  = HexToRgb "E7E6E6"   # lt2
$tcs.Item(5).RGB  = HexToRgb "5B9BD5"   # accent1
$tcs.Item(6).RGB  = HexToRgb "ED7D31"   # accent2
$tcs.Item(7).RGB  = HexToRgb "A5A5A5"   # accent3
$tcs.Item(8).RGB  = HexToRgb "FFC000"   # accent4
$tcs.Item(9).RGB  = HexToRgb "4472C4"   # accent5
$tcs.Item(10).RGB = HexToRgb "70AD47"   # accent6
$tcs.Item(11).RGB = HexToRgb "0563C1"   # hlink
$tcs.Item(12).RGB = HexToRgb "954F72"   # folHlink
